# Fruta / hortaliza, semanal
# Insert a new weekly record at row 456 (pushing the existing rows 456-473
# down to 457-474) in the "Feria Lagunitas de Puerto Montt - Naranja" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 456, shifting everything below
# (rows 456-473) down by one row.
$ws.Rows.Item(456).Insert()

# Populate the newly inserted row 456 with the new weekly price record.
$row = 456
$ws.Cells.Item($row, 1).Value  = 4
$ws.Cells.Item($row, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value  = "Los Lagos"
$ws.Cells.Item($row, 4).Value  = 44747
$ws.Cells.Item($row, 5).Value  = 10
$ws.Cells.Item($row, 6).Value  = "Fruta"
$ws.Cells.Item($row, 7).Value  = 100102
$ws.Cells.Item($row, 8).Value  = "Cítricos"
$ws.Cells.Item($row, 9).Value  = 100102005
$ws.Cells.Item($row, 10).Value = "Naranja"
$ws.Cells.Item($row, 11).Value = "Fukumoto"
$ws.Cells.Item($row, 12).Value = "Segunda"
$ws.Cells.Item($row, 13).Value = 800
$ws.Cells.Item($row, 14).Value = 9000
$ws.Cells.Item($row, 15).Value = 9000
$ws.Cells.Item($row, 16).Value = 9000
$ws.Cells.Item($row, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($row, 19).Value = 562
$ws.Cells.Item($row, 20).Value = 16
